# "Added New Mac-Address and Document Types"
#
# Appends 5 new device-master rows (id 3000176-3000180) to sheet1,
# following the exact same pattern as the preceding block of rows
# (Finger Print Scanner / IRIS Scanner / Web Camera / Document Scanner /
# Printer "31" -> "32"), and moves the current selection down to where
# the user was working (near row 113, full columns K:XFD selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ row=157; id=3000176; name="Finger Print Scanner 32"; mac="80-75-40-E8-CA-24"; serial="BS563Q2230824"; dspec=165 },
    @{ row=158; id=3000177; name="IRIS Scanner 32";         mac="0E-1A-14-4A-6D-3A"; serial="BS563Q2230825"; dspec=327 },
    @{ row=159; id=3000178; name="Web Camera 32";           mac="65-13-7F-0F-F7-53"; serial="BS563Q2230826"; dspec=736 },
    @{ row=160; id=3000179; name="Document Scanner 32";     mac="73-C4-DE-8E-C9-8D"; serial="BS563Q2230827"; dspec=801 },
    @{ row=161; id=3000180; name="Printer 32";              mac="EC-74-AB-E0-0F-38"; serial="BS563Q2230828"; dspec=920 }
)

# Columns: A id | B name | C mac_address | D serial_num | E ip_address (blank)
#        | F dspec_id | G lang_code | H is_active | I cr_by | J cr_dtimes
# Written column-by-column (all 5 rows for A, then all 5 for B, etc.) so
# that new shared-string entries land in the same order as the target
# workbook (all 5 names, then all 5 mac addresses, then all 5 serials).
foreach ($r in $newRows) {
    $ws.Cells.Item($r.row, 1).Value = $r.id
}
foreach ($r in $newRows) {
    $ws.Cells.Item($r.row, 2).Value = $r.name
}
foreach ($r in $newRows) {
    $ws.Cells.Item($r.row, 3).Value = $r.mac
}
foreach ($r in $newRows) {
    $ws.Cells.Item($r.row, 4).Value = $r.serial
}
foreach ($r in $newRows) {
    $ws.Cells.Item($r.row, 6).Value = $r.dspec
}
foreach ($r in $newRows) {
    $ws.Cells.Item($r.row, 7).Value = "eng"
}
foreach ($r in $newRows) {
    $ws.Cells.Item($r.row, 8).Value = $true
    # Matches the left-aligned style ("s=1") already used by every other
    # is_active cell in the column.
    $ws.Cells.Item($r.row, 8).HorizontalAlignment = -4131   # xlLeft
}
foreach ($r in $newRows) {
    $ws.Cells.Item($r.row, 9).Value = "superadmin"
}
foreach ($r in $newRows) {
    $ws.Cells.Item($r.row, 10).Value = "now()"
}

# Reflect the saved view state: scrolled so row 113 is at the top, with
# the full column range K:XFD selected (active cell K113).
$excel.ActiveWindow.ScrollRow = 113
$ws.Range("K1:XFD1048576").Select()
$excel.ActiveWindow.SetActiveCell($ws.Cells.Item(113, 11))
